$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# URL
$ws.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/identifier-type"

# Version
$ws.Range("B3").Value = "8.0.0"

# Date
$ws.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher
$ws.Range("B9").Value = "LinuxForHealth Team"

# Description
$ws.Range("B11").Value = "Extended set of Identifier type code for LinuxForHealth Common Data Model resources"
